$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance: TA Meeting on Wednesday Oct 6, 2021 (column D, row 10)
$ws.Range("D10").Value = 1

# Move active selection to the cell that was just edited
$ws.Range("D10").Select()
